# Updating LCOH & emissions charts for MN fact sheet
$wb = $excel.ActiveWorkbook

$baseline = $wb.Worksheets.Item("Baseline")
$pct80    = $wb.Worksheets.Item("80pct")
$pct100   = $wb.Worksheets.Item("100pct")
$readme   = $wb.Worksheets.Item("readme")

# --- Tidy up the readme note about the Baseline formula (set first so the
# --- shared-string table keeps the same append order as upstream).
$readme.Range("B5").Value = 'Baseline --> F29 "=F28*K29/1000"'

# --- Rename the "Waste heat pump" scenario columns to "Ambient heat pump" on
# --- all three data sheets, and shorten two sector row labels.
foreach ($ws in @($baseline, $pct80, $pct100)) {
    $ws.Range("I1").Value = "Energy efficiency + Ambient heat pump (worse case)"
    $ws.Range("J1").Value = "Energy efficiency + Ambient heat pump (best case)"
    $ws.Range("A3").Value = "Beet Sugar"
    $ws.Range("A4").Value = "Ethyl Alcohol"
}

# --- Updated emissions figures for the (renamed) Ambient heat pump columns.
$baseline.Range("I2").Value = 2198655
$baseline.Range("J2").Value = 2179010
$baseline.Range("I3").Value = 692856
$baseline.Range("J3").Value = 649484
$baseline.Range("I4").Value = 1427102
$baseline.Range("J4").Value = 1295341

$pct80.Range("I2").Value = 2001320
$pct80.Range("J2").Value = 1995368
$pct80.Range("I3").Value = 230921
$pct80.Range("J3").Value = 217778
$pct80.Range("I4").Value = 432455
$pct80.Range("J4").Value = 392528

# --- Move the active tab from 100pct to Baseline, and reset the selection
# --- on each data sheet to A2:A4.
$baseline.Activate()
$baseline.Range("A2:A4").Select()

$pct80.Activate()
$pct80.Range("A2:A4").Select()

$pct100.Activate()
$pct100.Range("A2:A4").Select()

$baseline.Activate()

# --- Match the author's resized/repositioned workbook window.
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 740
$win.Width = 29400
$win.Height = 18380
